$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 487-488, pushing the existing rows 487.. down to 489..
$ws.Rows("487:488").Insert()

# New row 487
$ws.Cells.Item(487, 1).Value = 1
$ws.Cells.Item(487, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(487, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(487, 4).Value = 45041
$ws.Cells.Item(487, 5).Value = 15
$ws.Cells.Item(487, 6).Value = 100114013
$ws.Cells.Item(487, 7).Value = "Zanahoria"
$ws.Cells.Item(487, 8).Value = "Sin especificar"
$ws.Cells.Item(487, 9).Value = "Primera"
$ws.Cells.Item(487, 10).Value = 60
$ws.Cells.Item(487, 11).Value = 14000
$ws.Cells.Item(487, 12).Value = 15000
$ws.Cells.Item(487, 13).Value = 14417
$ws.Cells.Item(487, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(487, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(487, 16).Value = 577
$ws.Cells.Item(487, 17).Value = 25
$ws.Cells.Item(487, 18).Value = "Hortaliza"

# New row 488
$ws.Cells.Item(488, 1).Value = 1
$ws.Cells.Item(488, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(488, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(488, 4).Value = 45041
$ws.Cells.Item(488, 5).Value = 15
$ws.Cells.Item(488, 6).Value = 100114013
$ws.Cells.Item(488, 7).Value = "Zanahoria"
$ws.Cells.Item(488, 8).Value = "Sin especificar"
$ws.Cells.Item(488, 9).Value = "Primera"
$ws.Cells.Item(488, 10).Value = 23
$ws.Cells.Item(488, 11).Value = 14000
$ws.Cells.Item(488, 12).Value = 15000
$ws.Cells.Item(488, 13).Value = 14435
$ws.Cells.Item(488, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(488, 15).Value = "Valle de Camiña"
$ws.Cells.Item(488, 16).Value = 577
$ws.Cells.Item(488, 17).Value = 25
$ws.Cells.Item(488, 18).Value = "Hortaliza"
